# Update "想去人数" (column F) counts on all sheets to reflect the latest
# scrape snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1613
$ws.Cells.Item(3, 6).Value = 850
$ws.Cells.Item(4, 6).Value = 256
$ws.Cells.Item(6, 6).Value = 1163
$ws.Cells.Item(7, 6).Value = 770
$ws.Cells.Item(8, 6).Value = 807
$ws.Cells.Item(9, 6).Value = 1485
$ws.Cells.Item(10, 6).Value = 298
$ws.Cells.Item(11, 6).Value = 1044
$ws.Cells.Item(13, 6).Value = 69
$ws.Cells.Item(16, 6).Value = 493
$ws.Cells.Item(17, 6).Value = 43
$ws.Cells.Item(18, 6).Value = 36
$ws.Cells.Item(22, 6).Value = 564
$ws.Cells.Item(23, 6).Value = 572
$ws.Cells.Item(24, 6).Value = 33
$ws.Cells.Item(26, 6).Value = 768

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 1014
$ws.Cells.Item(6, 6).Value = 15
$ws.Cells.Item(10, 6).Value = 86

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 260

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 260
$ws.Cells.Item(3, 6).Value = 1613
$ws.Cells.Item(5, 6).Value = 850
$ws.Cells.Item(6, 6).Value = 256
$ws.Cells.Item(7, 6).Value = 1014
$ws.Cells.Item(9, 6).Value = 1163
$ws.Cells.Item(10, 6).Value = 770
$ws.Cells.Item(11, 6).Value = 807
$ws.Cells.Item(12, 6).Value = 1485
$ws.Cells.Item(13, 6).Value = 298
$ws.Cells.Item(14, 6).Value = 1044
$ws.Cells.Item(16, 6).Value = 69
$ws.Cells.Item(19, 6).Value = 493
$ws.Cells.Item(20, 6).Value = 43
$ws.Cells.Item(21, 6).Value = 36
$ws.Cells.Item(27, 6).Value = 15
$ws.Cells.Item(30, 6).Value = 564
$ws.Cells.Item(31, 6).Value = 572
$ws.Cells.Item(32, 6).Value = 33
$ws.Cells.Item(34, 6).Value = 768
$ws.Cells.Item(39, 6).Value = 86
$ws.Cells.Item(40, 6).Value = 86

$wb.Save()
